# Update the dSF (column F) values to match the repulled/recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -1
    5  = 1
    9  = -3
    10 = -8
    13 = -2
    14 = -2
    22 = 7
    24 = -7
    25 = 4
    27 = -8
    32 = -5
    33 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
